$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (I1, J1) - match formatting of existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data values for columns I (I0) and J (IF), rows 2-12
$values = @(
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(6, 7),
    @(11, 11),
    @(7, 7),
    @(6, 7),
    @(8, 9),
    @(9, 9),
    @(6, 6),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
